$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 36666.332
$ws.Range("J3").Value = 36666.332
$ws.Range("L3").Value = 36666.332
$ws.Range("N3").Value = -36894.332

# Row 38
$ws.Range("H38").Value = 1985.8667
$ws.Range("I38").Value = 303.42856
$ws.Range("J38").Value = 3458
$ws.Range("K38").Value = 910.28568
$ws.Range("L38").Value = 10374
$ws.Range("M38").Value = -538.28568
$ws.Range("N38").Value = -11118

# Row 88
$ws.Range("H88").Value = 1847.8334
$ws.Range("I88").Value = 2178.5
$ws.Range("J88").Value = 1682.5
$ws.Range("K88").Value = 2178.5
$ws.Range("L88").Value = 1682.5
$ws.Range("M88").Value = -1772.5
$ws.Range("N88").Value = -2494.5

# Row 91
$ws.Range("H91").Value = 1847.8334
$ws.Range("I91").Value = 2178.5
$ws.Range("J91").Value = 1682.5
$ws.Range("K91").Value = 2178.5
$ws.Range("L91").Value = 1682.5
$ws.Range("M91").Value = -774.5
$ws.Range("N91").Value = -4490.5

# Row 94
$ws.Range("H94").Value = 3425.3333
$ws.Range("I94").Value = 3554.9092
$ws.Range("K94").Value = 3554.9092
$ws.Range("M94").Value = -3103.9092

# Row 102
$ws.Range("H102").Value = 36666.332
$ws.Range("J102").Value = 36666.332
$ws.Range("L102").Value = 36666.332
$ws.Range("N102").Value = -43156.332

# Row 116
$ws.Range("H116").Value = 6812
$ws.Range("J116").Value = 4124.5
$ws.Range("L116").Value = 4124.5
$ws.Range("N116").Value = -11008.5

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 11858.842
$ws.Range("J44").Value = 11858.842
$ws.Range("L44").Value = 11858.842
$ws.Range("N44").Value = -12834.842

# Row 55
$ws.Range("H55").Value = 99999
$ws.Range("J55").Value = 99999
$ws.Range("L55").Value = 99999
$ws.Range("N55").Value = -100629

# Row 101
$ws.Range("H101").Value = 20957.572
$ws.Range("J101").Value = 20957.572
$ws.Range("L101").Value = 20957.572
$ws.Range("N101").Value = -27447.572

# Row 110
$ws.Range("H110").Value = 1980.375
$ws.Range("I110").Value = 1398.1
$ws.Range("J110").Value = 2950.8333
$ws.Range("K110").Value = 1398.1
$ws.Range("L110").Value = 2950.8333
$ws.Range("M110").Value = 646.9000000000001
$ws.Range("N110").Value = -7040.8333

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3411.1875
$ws.Range("I86").Value = 1541.5
$ws.Range("J86").Value = 6527.3335
$ws.Range("K86").Value = 1541.5
$ws.Range("L86").Value = 6527.3335
$ws.Range("M86").Value = -418.5
$ws.Range("N86").Value = -8773.333500000001

# Row 89
$ws.Range("H89").Value = 3411.1875
$ws.Range("I89").Value = 1541.5
$ws.Range("J89").Value = 6527.3335
$ws.Range("K89").Value = 7707.5
$ws.Range("L89").Value = 32636.6675
$ws.Range("M89").Value = -2091.5
$ws.Range("N89").Value = -43868.6675

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6512.759
$ws.Range("I31").Value = 4773.5
$ws.Range("K31").Value = 4773.5
$ws.Range("M31").Value = -4478.5

# Row 34
$ws.Range("H34").Value = 6512.759
$ws.Range("I34").Value = 4773.5
$ws.Range("K34").Value = 4773.5
$ws.Range("M34").Value = -4571.5

# Row 43
$ws.Range("H43").Value = 50000
$ws.Range("J43").Value = 50000
$ws.Range("L43").Value = 50000
$ws.Range("N43").Value = -50368

# Row 58
$ws.Range("H58").Value = 3562.7778
$ws.Range("I58").Value = 1829.8334
$ws.Range("K58").Value = 1829.8334
$ws.Range("M58").Value = -1626.8334

# Row 93
$ws.Range("H93").Value = 6755.3335
$ws.Range("I93").Value = 6755.3335
$ws.Range("K93").Value = 6755.3335
$ws.Range("M93").Value = -4883.3335

# Row 99
$ws.Range("H99").Value = 4000
$ws.Range("I99").Value = 4000
$ws.Range("K99").Value = 4000
$ws.Range("M99").Value = -2502

# Row 101
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490

# Row 105
$ws.Range("H105").Value = 1821.6666
$ws.Range("I105").Value = 1479.75
$ws.Range("J105").Value = 2505.5
$ws.Range("K105").Value = 1479.75
$ws.Range("L105").Value = 2505.5
$ws.Range("M105").Value = 267.25
$ws.Range("N105").Value = -5999.5

# Row 126
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530

# Row 136
$ws.Range("H136").Value = 3562.7778
$ws.Range("I136").Value = 1829.8334
$ws.Range("K136").Value = 5489.5002
$ws.Range("M136").Value = -2939.5002

$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 4107.8125
$ws.Range("I18").Value = 2645
$ws.Range("K18").Value = 7935
$ws.Range("M18").Value = -7766

# Row 44
$ws.Range("H44").Value = 496.8
$ws.Range("I44").Value = 132.09091
$ws.Range("K44").Value = 396.27273
$ws.Range("M44").Value = 1.727269999999976

# Row 92
$ws.Range("H92").Value = 10001
$ws.Range("J92").Value = 10001.5
$ws.Range("L92").Value = 30004.5
$ws.Range("N92").Value = -32500.5

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 131
$ws.Range("H131").Value = 2990.889
$ws.Range("I131").Value = 2507.5
$ws.Range("J131").Value = 3377.6
$ws.Range("K131").Value = 7522.5
$ws.Range("L131").Value = 10132.8
$ws.Range("M131").Value = -2482.5
$ws.Range("N131").Value = -20212.8

# Row 134
$ws.Range("H134").Value = 3417.8
$ws.Range("I134").Value = 3417.8
$ws.Range("K134").Value = 10253.4
$ws.Range("M134").Value = -5183.400000000001

# Row 139
$ws.Range("H139").Value = 1956.9166
$ws.Range("I139").Value = 1156.4445
$ws.Range("K139").Value = 3469.3335
$ws.Range("M139").Value = 1670.6665

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 4908.778
$ws.Range("I43").Value = 2994.75
$ws.Range("J43").Value = 6440
$ws.Range("K43").Value = 2994.75
$ws.Range("L43").Value = 6440
$ws.Range("M43").Value = -2843.75
$ws.Range("N43").Value = -6742

# Row 44
$ws.Range("H44").Value = 4000
$ws.Range("J44").Value = 4000
$ws.Range("L44").Value = 4000
$ws.Range("N44").Value = -5192

# Row 70
$ws.Range("H70").Value = 1500
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# Row 73
$ws.Range("H73").Value = 1500
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# Row 94
$ws.Range("H94").Value = 18605.5
$ws.Range("J94").Value = 18605.5
$ws.Range("L94").Value = 18605.5
$ws.Range("N94").Value = -19957.5

# Row 95
$ws.Range("H95").Value = 32749.75
$ws.Range("J95").Value = 32749.75
$ws.Range("L95").Value = 32749.75
$ws.Range("N95").Value = -38241.75

# Row 101
$ws.Range("H101").Value = 27999.334
$ws.Range("J101").Value = 27999.334
$ws.Range("L101").Value = 27999.334
$ws.Range("N101").Value = -34489.334

# Row 107
$ws.Range("H107").Value = 209.75
$ws.Range("I107").Value = 209.75
$ws.Range("K107").Value = 209.75
$ws.Range("M107").Value = 1710.25

# Row 113
$ws.Range("H113").Value = 5351
$ws.Range("I113").Value = 3300.818
$ws.Range("K113").Value = 3300.818
$ws.Range("M113").Value = -1130.818

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 1176.75
$ws.Range("I55").Value = 1648
$ws.Range("K55").Value = 1648
$ws.Range("M55").Value = -1475

# Row 61
$ws.Range("H61").Value = 4472.5454
$ws.Range("I61").Value = 1840
$ws.Range("K61").Value = 1840
$ws.Range("M61").Value = -1638

# Row 93
$ws.Range("H93").Value = 1416.8667
$ws.Range("I93").Value = 1460.7858
$ws.Range("J93").Value = 802
$ws.Range("K93").Value = 1460.7858
$ws.Range("L93").Value = 802
$ws.Range("M93").Value = -212.7858000000001
$ws.Range("N93").Value = -3298

# Row 101
$ws.Range("H101").Value = 16936.857
$ws.Range("J101").Value = 16936.857
$ws.Range("L101").Value = 16936.857
$ws.Range("N101").Value = -23426.857

# Row 103
$ws.Range("H103").Value = 22100
$ws.Range("J103").Value = 22100
$ws.Range("L103").Value = 22100
$ws.Range("N103").Value = -24444

# Row 113
$ws.Range("H113").Value = 4472.5454
$ws.Range("I113").Value = 1840
$ws.Range("K113").Value = 1840
$ws.Range("M113").Value = 330

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 10389.429
$ws.Range("I51").Value = 9121
$ws.Range("J51").Value = 18000
$ws.Range("K51").Value = 9121
$ws.Range("L51").Value = 18000
$ws.Range("M51").Value = -8611
$ws.Range("N51").Value = -19020

# Row 113
$ws.Range("H113").Value = 916.5
$ws.Range("I113").Value = 1219.8
$ws.Range("J113").Value = 699.8570999999999
$ws.Range("K113").Value = 3659.4
$ws.Range("L113").Value = 2099.5713
$ws.Range("M113").Value = -1489.4
$ws.Range("N113").Value = -6439.5713

# Row 132
$ws.Range("H132").Value = 2749.5
$ws.Range("I132").Value = 2749.5
$ws.Range("K132").Value = 8248.5
$ws.Range("M132").Value = -5718.5
